$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
# row 18 (Leve Item ID 5471)
$ws.Range("H18").Value = 1454.5454
$ws.Range("I18").Value = 1333.3334
$ws.Range("K18").Value = 1333.3334
$ws.Range("M18").Value = -1049.3334
# row 28 (Leve Item ID 27772)
$ws.Range("H28").Value = 1834.5652
$ws.Range("I28").Value = 1904.3182
$ws.Range("K28").Value = 1904.3182
$ws.Range("M28").Value = -1419.3182
# row 33 (Leve Item ID 5512)
$ws.Range("H33").Value = 332.6154
$ws.Range("J33").Value = 291.83334
$ws.Range("L33").Value = 291.83334
$ws.Range("N33").Value = -749.83334
# row 40 (Leve Item ID 5505)
$ws.Range("H40").Value = 2015.3846
$ws.Range("I40").Value = 1911.1111
$ws.Range("K40").Value = 1911.1111
$ws.Range("M40").Value = -1736.1111
# row 53 (Leve Item ID 5479)
$ws.Range("H53").Value = 960.46155
$ws.Range("I53").Value = 1173.238
$ws.Range("K53").Value = 1173.238
$ws.Range("M53").Value = -536.2380000000001
# row 98 (Leve Item ID 36237)
$ws.Range("H98").Value = 1804.8197
$ws.Range("I98").Value = 1946.3704
$ws.Range("J98").Value = 712.8570999999999
$ws.Range("K98").Value = 1946.3704
$ws.Range("L98").Value = 712.8570999999999
$ws.Range("M98").Value = -448.3704
$ws.Range("N98").Value = -3708.8571
# row 105 (Leve Item ID 18668)
$ws.Range("H105").Value = 4671
$ws.Range("J105").Value = 4671
$ws.Range("L105").Value = 4671
$ws.Range("N105").Value = -11659
# row 113 (Leve Item ID 27775)
$ws.Range("H113").Value = 3935
$ws.Range("J113").Value = 3900
$ws.Range("L113").Value = 3900
$ws.Range("N113").Value = -10408
# row 116 (Leve Item ID 27778)
$ws.Range("H116").Value = 2774.6875
$ws.Range("I116").Value = 2038.9
$ws.Range("K116").Value = 2038.9
$ws.Range("M116").Value = 1403.1
# row 122 (Leve Item ID 36237)
$ws.Range("H122").Value = 1804.8197
$ws.Range("I122").Value = 1946.3704
$ws.Range("J122").Value = 712.8570999999999
$ws.Range("K122").Value = 5839.1112
$ws.Range("L122").Value = 2138.5713
$ws.Range("M122").Value = -3389.1112
$ws.Range("N122").Value = -7038.5713
# row 132 (Leve Item ID 44049)
$ws.Range("H132").Value = 10424205
$ws.Range("I132").Value = 11499995
$ws.Range("K132").Value = 34499985
$ws.Range("M132").Value = -34497455

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
# row 32 (Leve Item ID 44147)
$ws.Range("H32").Value = 10998.5
$ws.Range("I32").Value = 9330.634
$ws.Range("K32").Value = 9330.634
$ws.Range("M32").Value = -9043.634
# row 74 (Leve Item ID 44000)
$ws.Range("H74").Value = 1237.2142
$ws.Range("I74").Value = 1020.0909
$ws.Range("J74").Value = 2033.3334
$ws.Range("K74").Value = 1020.0909
$ws.Range("L74").Value = 2033.3334
$ws.Range("M74").Value = -146.0909
$ws.Range("N74").Value = -3781.3334
# row 77 (Leve Item ID 44000)
$ws.Range("H77").Value = 1237.2142
$ws.Range("I77").Value = 1020.0909
$ws.Range("J77").Value = 2033.3334
$ws.Range("K77").Value = 5100.4545
$ws.Range("L77").Value = 10166.667
$ws.Range("M77").Value = -732.4544999999998
$ws.Range("N77").Value = -18902.667
# row 132 (Leve Item ID 43997)
$ws.Range("H132").Value = 3479.7693
$ws.Range("I132").Value = 3658.7144
$ws.Range("K132").Value = 10976.1432
$ws.Range("M132").Value = -8446.143199999999
# row 141 (Leve Item ID 42483)
$ws.Range("H141").Value = 31583
$ws.Range("J141").Value = 31583
$ws.Range("L141").Value = 31583
$ws.Range("N141").Value = -41943

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
# row 80 (Leve Item ID 13747)
$ws.Range("H80").Value = 722.75
$ws.Range("I80").Value = 117.2
$ws.Range("J80").Value = 1155.2858
$ws.Range("K80").Value = 117.2
$ws.Range("L80").Value = 1155.2858
$ws.Range("M80").Value = 880.8
$ws.Range("N80").Value = -3151.2858
# row 83 (Leve Item ID 13747)
$ws.Range("H83").Value = 722.75
$ws.Range("I83").Value = 117.2
$ws.Range("J83").Value = 1155.2858
$ws.Range("K83").Value = 586
$ws.Range("L83").Value = 5776.429
$ws.Range("M83").Value = 4406
$ws.Range("N83").Value = -15760.429

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
# row 31 (Leve Item ID 44023)
$ws.Range("H31").Value = 1191.6809
$ws.Range("I31").Value = 984.18604
$ws.Range("J31").Value = 3422.25
$ws.Range("K31").Value = 984.18604
$ws.Range("L31").Value = 3422.25
$ws.Range("M31").Value = -689.18604
$ws.Range("N31").Value = -4012.25
# row 34 (Leve Item ID 44023)
$ws.Range("H34").Value = 1191.6809
$ws.Range("I34").Value = 984.18604
$ws.Range("J34").Value = 3422.25
$ws.Range("K34").Value = 984.18604
$ws.Range("L34").Value = 3422.25
$ws.Range("M34").Value = -782.18604
$ws.Range("N34").Value = -3826.25
# row 94 (Leve Item ID 32934)
$ws.Range("H94").Value = 1212
$ws.Range("I94").Value = 996.2857
$ws.Range("J94").Value = 1427.7142
$ws.Range("K94").Value = 996.2857
$ws.Range("L94").Value = 1427.7142
$ws.Range("M94").Value = -545.2857
$ws.Range("N94").Value = -2329.7142
# row 132 (Leve Item ID 44019)
$ws.Range("H132").Value = 7237.05
$ws.Range("I132").Value = 9227.154
$ws.Range("J132").Value = 3541.1428
$ws.Range("K132").Value = 27681.462
$ws.Range("L132").Value = 10623.4284
$ws.Range("M132").Value = -25151.462
$ws.Range("N132").Value = -15683.4284
# row 134 (Leve Item ID 44020)
$ws.Range("H134").Value = 20002112
$ws.Range("I134").Value = 2216.8333
$ws.Range("J134").Value = 71430420
$ws.Range("K134").Value = 6650.499899999999
$ws.Range("L134").Value = 214291260
$ws.Range("M134").Value = -4115.499899999999
$ws.Range("N134").Value = -214296330

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
# row 4 (Leve Item ID 4650)
$ws.Range("H4").Value = 188335.23
$ws.Range("I4").Value = 64264.215
$ws.Range("J4").Value = 767333.3
$ws.Range("K4").Value = 192792.645
$ws.Range("L4").Value = 2301999.9
$ws.Range("M4").Value = -192680.645
$ws.Range("N4").Value = -2302223.9
# row 18 (Leve Item ID 36056)
$ws.Range("H18").Value = 1500
$ws.Range("I18").Value = 1500
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 4500
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = -4331
$ws.Range("N18").ClearContents()
# row 107 (Leve Item ID 27838)
$ws.Range("H107").Value = 10390.4
$ws.Range("J107").Value = 16980.834
$ws.Range("L107").Value = 50942.50199999999
$ws.Range("N107").Value = -54782.50199999999
# row 130 (Leve Item ID 36058)
$ws.Range("H130").Value = 2330.9375
$ws.Range("I130").Value = 1000
$ws.Range("J130").Value = 2419.6667
$ws.Range("K130").Value = 3000
$ws.Range("L130").Value = 7259.000100000001
$ws.Range("M130").Value = 2020
$ws.Range("N130").Value = -17299.0001
# row 131 (Leve Item ID 36060)
$ws.Range("H131").Value = 21280074
$ws.Range("I131").Value = 200000290
$ws.Range("J131").Value = 3857.2856
$ws.Range("K131").Value = 600000870
$ws.Range("L131").Value = 11571.8568
$ws.Range("M131").Value = -599995830
$ws.Range("N131").Value = -21651.8568
# row 137 (Leve Item ID 44088)
$ws.Range("H137").Value = 25865296
$ws.Range("I137").Value = 107143860
$ws.Range("J137").Value = 3935.682
$ws.Range("K137").Value = 321431580
$ws.Range("L137").Value = 11807.046
$ws.Range("M137").Value = -321426480
$ws.Range("N137").Value = -22007.046

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
# row 124 (Leve Item ID 34247)
$ws.Range("H124").Value = 51074.75
$ws.Range("J124").Value = 51074.75
$ws.Range("L124").Value = 51074.75
$ws.Range("N124").Value = -60894.75
# row 126 (Leve Item ID 36184)
$ws.Range("H126").Value = 2484.2856
$ws.Range("I126").Value = 1866.6666
$ws.Range("K126").Value = 5599.9998
$ws.Range("M126").Value = -3129.9998
# row 132 (Leve Item ID 44008)
$ws.Range("H132").Value = 6126.125
$ws.Range("I132").Value = 8252.75
$ws.Range("J132").Value = 3999.5
$ws.Range("K132").Value = 24758.25
$ws.Range("L132").Value = 11998.5
$ws.Range("M132").Value = -22228.25
$ws.Range("N132").Value = -17058.5

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
# row 61 (Leve Item ID 27740)
$ws.Range("H61").Value = 1282
$ws.Range("I61").Value = 900
$ws.Range("J61").Value = 1377.5
$ws.Range("K61").Value = 900
$ws.Range("L61").Value = 1377.5
$ws.Range("M61").Value = -698
$ws.Range("N61").Value = -1781.5
# row 113 (Leve Item ID 27740)
$ws.Range("H113").Value = 1282
$ws.Range("I113").Value = 900
$ws.Range("J113").Value = 1377.5
$ws.Range("K113").Value = 900
$ws.Range("L113").Value = 1377.5
$ws.Range("M113").Value = 1270
$ws.Range("N113").Value = -5717.5
# row 122 (Leve Item ID 36247)
$ws.Range("H122").Value = 14716323
$ws.Range("I122").Value = 22738446
$ws.Range("J122").Value = 9099
$ws.Range("K122").Value = 68215338
$ws.Range("L122").Value = 27297
$ws.Range("M122").Value = -68212888
$ws.Range("N122").Value = -32197
# row 123 (Leve Item ID 35408)
$ws.Range("H123").Value = 40944
$ws.Range("J123").Value = 40944
$ws.Range("L123").Value = 40944
$ws.Range("N123").Value = -50744
# row 132 (Leve Item ID 44058)
$ws.Range("H132").Value = 76227.92999999999
$ws.Range("I132").Value = 22763.4
$ws.Range("J132").Value = 102960.2
$ws.Range("K132").Value = 68290.20000000001
$ws.Range("L132").Value = 308880.6
$ws.Range("M132").Value = -65760.20000000001
$ws.Range("N132").Value = -313940.6
# row 133 (Leve Item ID 41903)
$ws.Range("H133").Value = 39603.715
$ws.Range("J133").Value = 39603.715
$ws.Range("L133").Value = 39603.715
$ws.Range("N133").Value = -44663.715

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
# row 125 (Leve Item ID 34276)
$ws.Range("H125").Value = 39842.715
$ws.Range("J125").Value = 39842.715
$ws.Range("L125").Value = 39842.715
$ws.Range("N125").Value = -49682.715
# row 132 (Leve Item ID 44029)
$ws.Range("H132").Value = 4453.4
$ws.Range("I132").Value = 7934.3335
$ws.Range("K132").Value = 23803.0005
$ws.Range("M132").Value = -21273.0005
